# Revert to 2.1.1 files
#
# The original commit re-introduces a "Texas Notes" worksheet between the
# "Data" and "PPEIdtICEaT" sheets, containing a handful of free-text notes
# from the analyst, and leaves the cursor/selection on several sheets in
# new positions (artifacts of the author's last editing session).

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$data  = $wb.Worksheets.Item("Data")

# --- Insert the new "Texas Notes" sheet right after "Data" (i.e. right
# before "PPEIdtICEaT") -----------------------------------------------------
$texas = $wb.Worksheets.Add([System.Type]::Missing, $data)
$texas.Name = "Texas Notes"

# Re-fetch PPEIdtICEaT now that the sheet collection has changed, so the
# reference below isn't stale.
$ppe = $wb.Worksheets.Item("PPEIdtICEaT")

$texas.Range("A1").Value  = "This spreadsheet uses a very particular study. "
$texas.Range("A3").Value  = "It's done in Wisconsin where the authors use an educational seminar for builders"
$texas.Range("A4").Value  = "then they follow up with phone surveys to see what the builders actually implemented"
$texas.Range("A5").Value  = "then they use those results to try and estimate how much energy those builders decisions saved"
$texas.Range("A7").Value  = "It's all very niche and a bit subjective"
$texas.Range("A9").Value  = "That said, I did a quick literature search and didn't find anything that I thought"
$texas.Range("A10").Value = "would give us better or more Texas-specific numbers. "
$texas.Range("A12").Value = 'Since this is a "low" priority sheet, I will leave it alone.'

# --- Restore per-sheet cursor/selection state left behind by the author ---
$about.Activate()
$about.Range("B25").Select()

$data.Activate()
$data.Range("A14").Select()

$texas.Activate()
$texas.Range("C27").Select()

# PPEIdtICEaT is left as the active (tab-selected) sheet.
$ppe.Activate()
$ppe.Range("D20").Select()
